$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.56798742104897
$ws.Range("C2").Value = 5.5624014371473
$ws.Range("D2").Value = 8.904495809317133
$ws.Range("E2").Value = 10.67598311675275
$ws.Range("F2").Value = 43.5211027480406
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 14.91549994552507
$ws.Range("L2").Value = 10.11849393003563
$ws.Range("M2").Value = 16.99277315757973
$ws.Range("N2").Value = 24.89438636514307
$ws.Range("B3").Value = 18.42006533862364
$ws.Range("C3").Value = 5.399777418781572
$ws.Range("D3").Value = 8.913683858099295
$ws.Range("E3").Value = 10.69221696583141
$ws.Range("F3").Value = 43.39044630862689
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 14.81024531050813
$ws.Range("L3").Value = 10.12974653347591
$ws.Range("M3").Value = 16.98428539354111
$ws.Range("N3").Value = 24.9228411022896
$ws.Range("B4").Value = 18.3336389451732
$ws.Range("C4").Value = 5.296116411646329
$ws.Range("D4").Value = 8.919464365204654
$ws.Range("E4").Value = 10.70309845947532
$ws.Range("F4").Value = 43.31839245167113
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 14.74919779977409
$ws.Range("L4").Value = 10.13806042733651
$ws.Range("M4").Value = 16.98224569073491
$ws.Range("N4").Value = 24.94212439181826
$ws.Range("B5").Value = 18.29956018811472
$ws.Range("C5").Value = 5.25294532627352
$ws.Range("D5").Value = 8.921855124916888
$ws.Range("E5").Value = 10.70776292067688
$ws.Range("F5").Value = 43.29109947485035
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 14.72524368886243
$ws.Range("L5").Value = 10.14180184388382
$ws.Range("M5").Value = 16.98221387853328
$ws.Range("N5").Value = 24.95043806957895
$ws.Range("B6").Value = 18.29397130232191
$ws.Range("C6").Value = 5.245721643515979
$ws.Range("D6").Value = 8.922254237969245
$ws.Range("E6").Value = 10.70855136309512
$ws.Range("F6").Value = 43.28669289571207
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 14.72132253397171
$ws.Range("L6").Value = 10.14244445485612
$ws.Range("M6").Value = 16.98225692148549
$ws.Range("N6").Value = 24.95184606373734
$ws.Range("B7").Value = 18.33317468432157
$ws.Range("C7").Value = 5.295537909298368
$ws.Range("D7").Value = 8.91949646523171
$ws.Range("E7").Value = 10.70316043360817
$ws.Range("F7").Value = 43.31801597110373
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 14.74887097868088
$ws.Range("L7").Value = 10.13810945410836
$ws.Range("M7").Value = 16.98224202297648
$ws.Range("N7").Value = 24.94223466848531
$ws.Range("B8").Value = 18.51609248382137
$ws.Range("C8").Value = 5.507136106534033
$ws.Range("D8").Value = 8.907635120371008
$ws.Range("E8").Value = 10.68139110247533
$ws.Range("F8").Value = 43.47436652508994
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 14.878479700876
$ws.Range("L8").Value = 10.12208240659654
$ws.Range("M8").Value = 16.98918946595225
$ws.Range("N8").Value = 24.90382145390389
$ws.Range("B9").Value = 18.90806732960904
$ws.Range("C9").Value = 5.890700165091443
$ws.Range("D9").Value = 8.885468287320791
$ws.Range("E9").Value = 10.64593616881411
$ws.Range("F9").Value = 43.84503576073716
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 15.15994497189854
$ws.Range("L9").Value = 10.10178950066171
$ws.Range("M9").Value = 17.02788550263548
$ws.Range("N9").Value = 24.84287339385068
$ws.Range("B10").Value = 19.21403278915723
$ws.Range("C10").Value = 6.152014953496765
$ws.Range("D10").Value = 8.869834362736178
$ws.Range("E10").Value = 10.62427594959175
$ws.Range("F10").Value = 44.15528379456631
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 15.38181409263101
$ws.Range("L10").Value = 10.09365403966963
$ws.Range("M10").Value = 17.07145004900973
$ws.Range("N10").Value = 24.80686640754907
$ws.Range("B11").Value = 19.3566091070741
$ws.Range("C11").Value = 6.266205960724743
$ws.Range("D11").Value = 8.862860524544434
$ws.Range("E11").Value = 10.6153704080987
$ws.Range("F11").Value = 44.30437325575988
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 15.48566327898226
$ws.Range("L11").Value = 10.0914195379164
$ws.Range("M11").Value = 17.09451428978783
$ws.Range("N11").Value = 24.7923916567446
$ws.Range("B12").Value = 19.41104158599904
$ws.Range("C12").Value = 6.308756699414012
$ws.Range("D12").Value = 8.860239344278325
$ws.Range("E12").Value = 10.6121340236686
$ws.Range("F12").Value = 44.36194650220705
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 15.52537665974701
$ws.Range("L12").Value = 10.09078376085401
$ws.Range("M12").Value = 17.10371074812149
$ws.Range("N12").Value = 24.78718446527027
$ws.Range("B13").Value = 19.39929965718623
$ws.Range("C13").Value = 6.299623623678801
$ws.Range("D13").Value = 8.860802991098321
$ws.Range("E13").Value = 10.61282499598238
$ws.Range("F13").Value = 44.34949788648812
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 15.51680693345432
$ws.Range("L13").Value = 10.09091133844325
$ws.Range("M13").Value = 17.10170962835001
$ws.Range("N13").Value = 24.7882937349869
$ws.Range("B14").Value = 19.36107870704398
$ws.Range("C14").Value = 6.269720576118456
$ws.Range("D14").Value = 8.862644485577276
$ws.Range("E14").Value = 10.61510142611234
$ws.Range("F14").Value = 44.30908763436953
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 15.48892293275063
$ws.Range("L14").Value = 10.0913630188037
$ws.Range("M14").Value = 17.09526164438142
$ws.Range("N14").Value = 24.79195776414608
$ws.Range("B15").Value = 19.33772345572381
$ws.Range("C15").Value = 6.251313631087522
$ws.Range("D15").Value = 8.863775008619848
$ws.Range("E15").Value = 10.61651349992701
$ws.Range("F15").Value = 44.2844797092956
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 15.47189274115117
$ws.Range("L15").Value = 10.09166706841935
$ws.Range("M15").Value = 17.09137216781322
$ws.Range("N15").Value = 24.79423778560092
$ws.Range("B16").Value = 19.2047794359331
$ws.Range("C16").Value = 6.144456529471659
$ws.Range("D16").Value = 8.870292881590077
$ws.Range("E16").Value = 10.62487699069295
$ws.Range("F16").Value = 44.14569821113681
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 15.37508334511029
$ws.Range("L16").Value = 10.09382953692733
$ws.Range("M16").Value = 17.07000771589241
$ws.Range("N16").Value = 24.80785071271188
$ws.Range("B17").Value = 19.1240588165282
$ws.Range("C17").Value = 6.077690702880173
$ws.Range("D17").Value = 8.874326617571558
$ws.Range("E17").Value = 10.63025023355324
$ws.Range("F17").Value = 44.06257881185354
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 15.31641926179651
$ws.Range("L17").Value = 10.09553138990165
$ws.Range("M17").Value = 17.05772989457792
$ws.Range("N17").Value = 24.81668982590658
$ws.Range("B18").Value = 19.07795242391676
$ws.Range("C18").Value = 6.038848950275729
$ws.Range("D18").Value = 8.876659727571516
$ws.Range("E18").Value = 10.63343001284723
$ws.Range("F18").Value = 44.01552137454548
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 15.28295401885324
$ws.Range("L18").Value = 10.09664830434403
$ws.Range("M18").Value = 17.05097387928074
$ws.Range("N18").Value = 24.82195315023647
$ws.Range("B19").Value = 19.06239827272178
$ws.Range("C19").Value = 6.025622831771217
$ws.Range("D19").Value = 8.877451919377515
$ws.Range("E19").Value = 10.63452196754773
$ws.Range("F19").Value = 43.99971827508005
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 15.2716717621364
$ws.Range("L19").Value = 10.09705019679844
$ws.Range("M19").Value = 17.04873906510531
$ws.Range("N19").Value = 24.82376601184054
$ws.Range("B20").Value = 19.13261867052931
$ws.Range("C20").Value = 6.084843692201631
$ws.Range("D20").Value = 8.87389587429546
$ws.Range("E20").Value = 10.62966900996399
$ws.Range("F20").Value = 44.07134952088413
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 15.32263573644945
$ws.Range("L20").Value = 10.0953359400163
$ws.Range("M20").Value = 17.05900526390976
$ws.Range("N20").Value = 24.81573032936401
$ws.Range("B21").Value = 19.37229350019715
$ws.Range("C21").Value = 6.278522707415924
$ws.Range("D21").Value = 8.862103062009147
$ws.Range("E21").Value = 10.61442909647135
$ws.Range("F21").Value = 44.32092703889838
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 15.49710286840505
$ws.Range("L21").Value = 10.09122464382411
$ws.Range("M21").Value = 17.09714305713102
$ws.Range("N21").Value = 24.79087411032331
$ws.Range("B22").Value = 19.53148915514727
$ws.Range("C22").Value = 6.401069196810713
$ws.Range("D22").Value = 8.854510283510553
$ws.Range("E22").Value = 10.60526119724492
$ws.Range("F22").Value = 44.49053293958847
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 15.61337326886032
$ws.Range("L22").Value = 10.08976366244033
$ws.Range("M22").Value = 17.12476199022953
$ws.Range("N22").Value = 24.77622675897767
$ws.Range("B23").Value = 19.44630496032119
$ws.Range("C23").Value = 6.336038225319316
$ws.Range("D23").Value = 8.858552280988816
$ws.Range("E23").Value = 10.61008189885095
$ws.Range("F23").Value = 44.3994266338031
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 15.55112275419732
$ws.Range("L23").Value = 10.09043141321245
$ws.Range("M23").Value = 17.10977630994201
$ws.Range("N23").Value = 24.78389809742718
$ws.Range("B24").Value = 19.12874781959426
$ws.Range("C24").Value = 6.081611248456417
$ws.Range("D24").Value = 8.874090569502698
$ws.Range("E24").Value = 10.6299314988287
$ws.Range("F24").Value = 44.06738201040583
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 15.31982445123871
$ws.Range("L24").Value = 10.09542387147957
$ws.Range("M24").Value = 17.05842772634191
$ws.Range("N24").Value = 24.81616355206206
$ws.Range("B25").Value = 18.79869717418556
$ws.Range("C25").Value = 5.790445355787173
$ws.Range("D25").Value = 8.89134945460375
$ws.Range("E25").Value = 10.65475540528343
$ws.Range("F25").Value = 43.73801068768994
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 15.08103502444431
$ws.Range("L25").Value = 10.10608839980892
$ws.Range("M25").Value = 17.01474682229505
$ws.Range("N25").Value = 24.85782142415323
